$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values are NOT numeric-looking: safe to assign directly as strings.
$plainUpdates = @(
    @{ Cell = 'D2'; Value = '25.722.64' }
    @{ Cell = 'E2'; Value = '  -2.83%  ' }
    @{ Cell = 'D3'; Value = '1.745.42' }
    @{ Cell = 'E3'; Value = '  -5.20%  ' }
    @{ Cell = 'E4'; Value = '  +0.17%  ' }
    @{ Cell = 'E5'; Value = '  -8.81%  ' }
    @{ Cell = 'E6'; Value = '  +0.13%  ' }
    @{ Cell = 'E7'; Value = '  -6.07%  ' }
    @{ Cell = 'E8'; Value = '  -6.68%  ' }
    @{ Cell = 'E9'; Value = '  -13.54%  ' }
    @{ Cell = 'E10'; Value = '  -10.89%  ' }
    @{ Cell = 'D11'; Value = '1.748.79' }
    @{ Cell = 'E11'; Value = '  -5.16%  ' }
    @{ Cell = 'E12'; Value = '  -8.31%  ' }
    @{ Cell = 'E13'; Value = '  -15.29%  ' }
    @{ Cell = 'E14'; Value = '  -10.18%  ' }
    @{ Cell = 'E15'; Value = '  -20.02%  ' }
    @{ Cell = 'E16'; Value = '  -14.45%  ' }
    @{ Cell = 'E17'; Value = '  +0.11%  ' }
    @{ Cell = 'E18'; Value = '  +0.15%  ' }
    @{ Cell = 'D19'; Value = '25.774.33' }
    @{ Cell = 'E19'; Value = '  -2.75%  ' }
    @{ Cell = 'E20'; Value = '  -16.86%  ' }
    @{ Cell = 'E21'; Value = '  -14.52%  ' }
    @{ Cell = 'D22'; Value = '1.971.59' }
    @{ Cell = 'E22'; Value = '  -5.52%  ' }
    @{ Cell = 'E23'; Value = '  -11.72%  ' }
    @{ Cell = 'E24'; Value = '  -12.97%  ' }
    @{ Cell = 'E25'; Value = '  -14.78%  ' }
    @{ Cell = 'E26'; Value = '  -3.40%  ' }
    @{ Cell = 'E27'; Value = '  -8.76%  ' }
    @{ Cell = 'E28'; Value = '  -18.07%  ' }
    @{ Cell = 'E29'; Value = '  -12.10%  ' }
    @{ Cell = 'E30'; Value = '  -6.84%  ' }
    @{ Cell = 'E31'; Value = '  -11.66%  ' }
    @{ Cell = 'E32'; Value = '  -7.91%  ' }
    @{ Cell = 'E33'; Value = '  -15.08%  ' }
    @{ Cell = 'E34'; Value = '  -6.42%  ' }
    @{ Cell = 'E35'; Value = '  +0.05%  ' }
    @{ Cell = 'E36'; Value = '  -10.08%  ' }
    @{ Cell = 'E37'; Value = '  -13.91%  ' }
    @{ Cell = 'E38'; Value = '  -17.09%  ' }
    @{ Cell = 'E39'; Value = '  -14.58%  ' }
    @{ Cell = 'E40'; Value = '  -10.25%  ' }
    @{ Cell = 'E41'; Value = '  -3.23%  ' }
    @{ Cell = 'E42'; Value = '  +0.10%  ' }
    @{ Cell = 'E43'; Value = '  -17.05%  ' }
    @{ Cell = 'E44'; Value = '  -12.67%  ' }
    @{ Cell = 'E45'; Value = '  -20.78%  ' }
    @{ Cell = 'E46'; Value = '  -19.82%  ' }
    @{ Cell = 'E47'; Value = '  -8.19%  ' }
    @{ Cell = 'E48'; Value = '  -10.44%  ' }
    @{ Cell = 'B49'; Value = 'Elrond' }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld' }
    @{ Cell = 'E49'; Value = '  -14.20%  ' }
    @{ Cell = 'B50'; Value = 'Aptos' }
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt' }
    @{ Cell = 'E50'; Value = '  -21.13%  ' }
    @{ Cell = 'E51'; Value = '  -13.48%  ' })

foreach ($u in $plainUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}

# Cells whose new text values look like plain numbers (e.g. "1.004", "41.72").
# Excel would normally auto-convert these to actual numeric values, losing the
# original text formatting (e.g. trailing zeros such as "1.004" vs 1.004,
# "11.60" vs 11.6). Force them to remain text by marking the cell as Text
# before assigning the value, then clear the (temporary) number-format again
# so no stray cell style is left behind - only the text content changes.
$numericLooksLikeUpdates = @(
    @{ Cell = 'D4'; Value = '1.004' }
    @{ Cell = 'D5'; Value = '238.26' }
    @{ Cell = 'D7'; Value = '0.5029' }
    @{ Cell = 'D8'; Value = '41.72' }
    @{ Cell = 'D9'; Value = '0.2622' }
    @{ Cell = 'D10'; Value = '0.06137' }
    @{ Cell = 'D12'; Value = '0.06955' }
    @{ Cell = 'D13'; Value = '15.21' }
    @{ Cell = 'D14'; Value = '4.478' }
    @{ Cell = 'D15'; Value = '0.5893' }
    @{ Cell = 'D16'; Value = '76.74' }
    @{ Cell = 'D17'; Value = '1.004' }
    @{ Cell = 'D18'; Value = '1.003' }
    @{ Cell = 'D20'; Value = '11.60' }
    @{ Cell = 'D21'; Value = '0.000006781' }
    @{ Cell = 'D23'; Value = '4.059' }
    @{ Cell = 'D24'; Value = '8.093' }
    @{ Cell = 'D25'; Value = '5.096' }
    @{ Cell = 'D26'; Value = '138.37' }
    @{ Cell = 'D27'; Value = '1.538' }
    @{ Cell = 'D28'; Value = '1.813' }
    @{ Cell = 'D29'; Value = '14.92' }
    @{ Cell = 'D30'; Value = '103.18' }
    @{ Cell = 'D31'; Value = '3.761' }
    @{ Cell = 'D32'; Value = '0.08109' }
    @{ Cell = 'D33'; Value = '3.446' }
    @{ Cell = 'D34'; Value = '0.04493' }
    @{ Cell = 'D36'; Value = '2.634' }
    @{ Cell = 'D37'; Value = '0.9752' }
    @{ Cell = 'D38'; Value = '0.6021' }
    @{ Cell = 'D39'; Value = '2.652' }
    @{ Cell = 'D40'; Value = '0.01541' }
    @{ Cell = 'D41'; Value = '104.44' }
    @{ Cell = 'D42'; Value = '1.002' }
    @{ Cell = 'D43'; Value = '1.903' }
    @{ Cell = 'D44'; Value = '5.129' }
    @{ Cell = 'D45'; Value = '0.3773' }
    @{ Cell = 'D46'; Value = '0.7278' }
    @{ Cell = 'D47'; Value = '0.05324' }
    @{ Cell = 'D48'; Value = '0.1105' }
    @{ Cell = 'D49'; Value = '29.97' }
    @{ Cell = 'D50'; Value = '5.876' }
    @{ Cell = 'D51'; Value = '52.24' })

foreach ($u in $numericLooksLikeUpdates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.ClearFormats()
}
